$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "308.71"
    "E2" = "-2.00%"
    "D3" = "38.06"
    "E3" = "-3.48%"
    "E4" = "-1.66%"
    "D5" = "0.07911"
    "E5" = "-3.42%"
    "D6" = "2.024"
    "E6" = "2.25%"
    "D7" = "4.419"
    "E7" = "4.32%"
    "D8" = "8.258"
    "E8" = "1.22%"
    "D9" = "3.212"
    "E9" = "-2.56%"
    "D10" = "0.9294"
    "E10" = "0.23%"
    "D11" = "0.1280"
    "E11" = "-7.65%"
    "D12" = "0.1891"
    "E12" = "-3.86%"
    "D13" = "0.08751"
    "E13" = "-3.17%"
    "D14" = "0.03453"
    "E14" = "-1.55%"
    "D15" = "0.09747"
    "E15" = "-0.71%"
    "D16" = "0.001393"
    "E16" = "0.00%"
    "D17" = "0.006098"
    "E17" = "-0.08%"
    "D18" = "3.590"
    "E18" = "-2.34%"
    "D19" = "0.3443"
    "E19" = "-0.56%"
    "D20" = "0.1296"
    "E20" = "-4.08%"
    "D21" = "5.014"
    "E21" = "8.07%"
    "D22" = "0.2518"
    "E22" = "4.07%"
    "D23" = "0.04340"
    "E23" = "-0.84%"
    "D24" = "0.001219"
    "E24" = "-1.93%"
    "D25" = "0.004627"
    "E25" = "-3.62%"
    "E26" = "177.00%"
    "D39" = "0.02268"
    "E39" = "4.94%"
    "D40" = "0.05077"
    "E40" = "-2.63%"
    "D41" = "0.007503"
    "E41" = "1.25%"
    "D42" = "0.009907"
    "E42" = "0.95%"
    "D43" = "0.1367"
    "E43" = "-0.51%"
    "D44" = "0.002095"
    "E44" = "-1.42%"
    "D45" = "0.008833"
    "E45" = "-10.55%"
    "D46" = "0.00006496"
    "E46" = "1.70%"
    "D47" = "0.00000000752"
    "E47" = "0.40%"
    "D48" = "0.003005"
    "E48" = "8.74%"
    "D49" = "0.001203"
    "E49" = "20.47%"
    "D50" = "0.00002105"
    "E50" = "0.40%"
    "D51" = "0.0002005"
    "E51" = "0.40%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}

